# The edit cyclically shifts the species-record data in rows 2-5 up by one
# row (row2 -> row5, row3 -> row2, row4 -> row3, row5 -> row4), while
# leaving the shared/common columns (C, I, P, S, T, U, V, W, Y, Z, AA, AB,
# AD, AE, AG, AT, AW, AX, AY) untouched since they are identical across
# these rows. Column AF only ever holds a single empty marker cell that
# moves from row 5 to row 4 as part of the same shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that vary per-record.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

$data = @{}
foreach ($row in 2..5) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value()
    }
    $data[$row] = $rowData
}

# before-row -> after-row mapping (content of $src ends up in $dst)
$mapping = @{ 2 = 5; 3 = 2; 4 = 3; 5 = 4 }

foreach ($src in $mapping.Keys) {
    $dst = $mapping[$src]
    foreach ($col in $cols) {
        $ws.Range("$col$dst").Value = $data[$src][$col]
    }
}

# Move the lone empty AF marker cell from row 5 to row 4 (Cut leaves the
# source cell completely empty/removed and creates a blank cell at the
# destination, matching the original empty-inlineStr marker semantics).
$ws.Range("AF5").Cut($ws.Range("AF4"))
